$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.785.96"
$ws.Range("E2").Value = "  -1.82%  "
$ws.Range("D3").Value = "3.483.63"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.29"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "129.57"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.98%  "
$ws.Range("D7").Value = "3.481.73"
$ws.Range("E7").Value = "  -1.09%  "
$ws.Range("E8").Value = "  -0.11%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.481"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.45%  "
$ws.Range("E10").Value = "  -1.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.07"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.18%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.375"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.21%  "
$ws.Range("D13").Value = "4.072.24"
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "27.13"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.05%  "
$ws.Range("E15").Value = "  +1.01%  "
$ws.Range("D16").Value = "3.498.61"
$ws.Range("E16").Value = "  -0.71%  "
$ws.Range("E17").Value = "  -3.15%  "
$ws.Range("D18").Value = "63.797.85"
$ws.Range("E18").Value = "  -1.91%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "9.77"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.81%  "
$ws.Range("E20").Value = "  -2.35%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "5.60"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.94%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "378.86"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.45%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.569"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -2.00%  "
$ws.Range("D24").Value = "3.617.99"
$ws.Range("E24").Value = "  -1.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "72.92"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.01%  "
$ws.Range("E26").Value = "  +0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0000114"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.84%  "
$ws.Range("E28").Value = "  -1.31%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.17%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -3.51%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.18"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.20%  "
$ws.Range("E32").Value = "  -2.67%  "
$ws.Range("D33").Value = "3.488.05"
$ws.Range("E33").Value = "  -1.32%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("E35").Value = "  -3.71%  "
$ws.Range("E36").Value = "  -1.26%  "
$ws.Range("E37").Value = "  -0.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.85"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("E39").Value = "  -2.28%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "159.49"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.34%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0788"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.49%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.807"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "25.97"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.43%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.75"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.19"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -5.37%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.92%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.41%  "
$ws.Range("D49").Value = "2.412.70"
$ws.Range("E49").Value = "  -0.05%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.79"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("E51").Value = "  -1.88%  "
